$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New enum block: "EnumCardProperty" with items Virtual/虚幻 and Flash/闪回
$ws.Range("B22").Value = "EnumCardProperty"
$ws.Range("D22").Value = $true
$ws.Range("G22").Value = "Virtual"
$ws.Range("J22").Value = "虚幻"

$ws.Range("G23").Value = "Flash"
$ws.Range("J23").Value = "闪回"

# Update the saved selection to match the authored state
$ws.Range("K16").Select()
